$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# =========================================================================
# 1) Update existing "debug stats" (Min/Max Conv percentages) in columns R/T
# =========================================================================
$ws.Range("R2").Value = 0.03
$ws.Range("T2").Value = 0.05

$ws.Range("R3").Value = -0.06
$ws.Range("T3").Value = -0.03

$ws.Range("R5").Value = -0.06

$ws.Range("R6").Value = -0.06

$ws.Range("R7").Value = -0.04

$ws.Range("R11").Value = -0.1
$ws.Range("T11").Value = -0.05

$ws.Range("T21").Value = -0.01

$ws.Range("R25").Value = -0.05

$ws.Range("R26").Value = -0.07

$ws.Range("R27").Value = -0.1

$ws.Range("T31").Value = 0.03

$ws.Range("R32").Value = -0.1

$ws.Range("R33").Value = -0.1

$ws.Range("R35").Value = -0.02

$ws.Range("R50").Value = -0.1

$ws.Range("R51").Value = -0.06

$ws.Range("R52").Value = -0.06

$ws.Range("R53").Value = -0.06

$ws.Range("R54").Value = 0.02
$ws.Range("T54").Value = 0.03

$ws.Range("T57").Value = 0.03

$ws.Range("R58").Value = -0.06

$ws.Range("T59").Value = 0.03

$ws.Range("R60").Value = -0.06

$ws.Range("T61").Value = 0.03

$ws.Range("T62").Value = 0.03

$ws.Range("T63").Value = 0.03

$ws.Range("T64").Value = 0.03

$ws.Range("T65").Value = 0.03

$ws.Range("T66").Value = 0.03

$ws.Range("T67").Value = 0.03

$ws.Range("T68").Value = 0.03

$ws.Range("T69").Value = 0.03

$ws.Range("T70").Value = 0.03

$ws.Range("T77").Value = 0.03

$ws.Range("T78").Value = 0.03

$ws.Range("R81").Value = -0.1

# =========================================================================
# 2) Add the new "endgame" news entries as rows 82-86
# =========================================================================
$styledCols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA")

# --- Row 82 (News ID 81) ---
$r = 82
foreach ($col in $styledCols) { $ws.Range("$col$r").Style = "Good" }
$ws.Range("A$r").Value = "new News( "
$ws.Range("B$r").Value = 81
$ws.Range("C$r").Value = ","
$ws.Range("D$r").Value = "`"ExPresidente Palometa, exiliado en Murcielandia`""
$ws.Range("E$r").Value = ","
$ws.Range("F$r").Value = "`"El ExPresidente fue recibido por autoridades de Murcielandia, en donde se encuentra exiliado.`""
$ws.Range("G$r").Value = ","
$ws.Range("H$r").Value = "`"`""
$ws.Range("I$r").Value = ","
$ws.Range("J$r").Value = 0.03
$ws.Range("K$r").Value = ","
$ws.Range("L$r").Value = 0.06
$ws.Range("M$r").Value = ","
$ws.Range("N$r").Value = -0.03
$ws.Range("O$r").Value = ","
$ws.Range("P$r").Value = -0.02
$ws.Range("Q$r").Value = ","
$ws.Range("R$r").Value = -0.1
$ws.Range("S$r").Value = ","
$ws.Range("T$r").Value = -0.02
$ws.Range("U$r").Value = ","
$ws.Range("V$r").Value = 0
$ws.Range("W$r").Value = ", new List<int>() {"
$ws.Range("Y$r").Value = "}, new List<int>(){"
$ws.Range("AA$r").Value = "}"
$ws.Range("AB$r").Value = "),"

# --- Row 83 (News ID 82) ---
$r = 83
foreach ($col in $styledCols) { $ws.Range("$col$r").Style = "Good" }
$ws.Range("A$r").Value = "new News( "
$ws.Range("B$r").Value = 82
$ws.Range("C$r").Value = ","
$ws.Range("D$r").Value = "`"Devastadoras declaraciones de ExPresidente Palometa`""
$ws.Range("E$r").Value = ","
$ws.Range("F$r").Value = "`"Dijo: Me taparon los ojos, y me obligaron a dejar el palacio presidencial. Es un golpe de lo mas ruin.`""
$ws.Range("G$r").Value = ","
$ws.Range("H$r").Value = "`"`""
$ws.Range("I$r").Value = ","
$ws.Range("J$r").Value = 0.03
$ws.Range("K$r").Value = ","
$ws.Range("L$r").Value = 0.06
$ws.Range("M$r").Value = ","
$ws.Range("N$r").Value = -0.03
$ws.Range("O$r").Value = ","
$ws.Range("P$r").Value = -0.02
$ws.Range("Q$r").Value = ","
$ws.Range("R$r").Value = -0.1
$ws.Range("S$r").Value = ","
$ws.Range("T$r").Value = -0.02
$ws.Range("U$r").Value = ","
$ws.Range("V$r").Value = 81
$ws.Range("W$r").Value = ", new List<int>() {"
$ws.Range("Y$r").Value = "}, new List<int>(){"
$ws.Range("AA$r").Value = "}"
$ws.Range("AB$r").Value = "),"

# --- Row 84 (News ID 83) ---
$r = 84
foreach ($col in $styledCols) { $ws.Range("$col$r").Style = "Good" }
$ws.Range("A$r").Value = "new News( "
$ws.Range("B$r").Value = 83
$ws.Range("C$r").Value = ","
$ws.Range("D$r").Value = "`"ExPresidente Palometa se ha suicidado`""
$ws.Range("E$r").Value = ","
$ws.Range("F$r").Value = "`"Fue encontrado en dudosas circunstancias. El peritaje inicial no es concluyente.`""
$ws.Range("G$r").Value = ","
$ws.Range("H$r").Value = "`"`""
$ws.Range("I$r").Value = ","
$ws.Range("J$r").Value = 0.03
$ws.Range("K$r").Value = ","
$ws.Range("L$r").Value = 0.06
$ws.Range("M$r").Value = ","
$ws.Range("N$r").Value = -0.03
$ws.Range("O$r").Value = ","
$ws.Range("P$r").Value = -0.02
$ws.Range("Q$r").Value = ","
$ws.Range("R$r").Value = -0.1
$ws.Range("S$r").Value = ","
$ws.Range("T$r").Value = -0.02
$ws.Range("U$r").Value = ","
$ws.Range("V$r").Value = 82
$ws.Range("W$r").Value = ", new List<int>() {"
$ws.Range("Y$r").Value = "}, new List<int>(){"
$ws.Range("AA$r").Value = "}"
$ws.Range("AB$r").Value = "),"

# --- Row 85 (News ID 84) ---
$r = 85
foreach ($col in $styledCols) { $ws.Range("$col$r").Style = "Good" }
$ws.Range("A$r").Value = "new News( "
$ws.Range("B$r").Value = 84
$ws.Range("C$r").Value = ","
$ws.Range("D$r").Value = "`"Autopsia de Palometa`""
$ws.Range("E$r").Value = ","
$ws.Range("F$r").Value = "`"Los resultados indican que se trataría de un magnicidio, servicios de inteligencia de Albatros, sospechados.`""
$ws.Range("G$r").Value = ","
$ws.Range("H$r").Value = "`"`""
$ws.Range("I$r").Value = ","
$ws.Range("J$r").Value = 0.03
$ws.Range("K$r").Value = ","
$ws.Range("L$r").Value = 0.06
$ws.Range("M$r").Value = ","
$ws.Range("N$r").Value = -0.03
$ws.Range("O$r").Value = ","
$ws.Range("P$r").Value = -0.02
$ws.Range("Q$r").Value = ","
$ws.Range("R$r").Value = -0.1
$ws.Range("S$r").Value = ","
$ws.Range("T$r").Value = -0.02
$ws.Range("U$r").Value = ","
$ws.Range("V$r").Value = 83
$ws.Range("W$r").Value = ", new List<int>() {"
$ws.Range("Y$r").Value = "}, new List<int>(){"
$ws.Range("AA$r").Value = "}"
$ws.Range("AB$r").Value = "),"

# --- Row 86 (News ID 85) ---
$r = 86
foreach ($col in $styledCols) { $ws.Range("$col$r").Style = "Bad" }
$ws.Range("A$r").Value = "new News( "
$ws.Range("B$r").Value = 85
$ws.Range("C$r").Value = ","
$ws.Range("D$r").Value = "`"La autopsia de Palometa es un circo`""
$ws.Range("E$r").Value = ","
$ws.Range("F$r").Value = "`"Expertos internacionales indican que la autopsia es un circo, que indiscutiblemente fue suicidio por la situación en la que dejó el país.`""
$ws.Range("G$r").Value = ","
$ws.Range("H$r").Value = "`"`""
$ws.Range("I$r").Value = ","
$ws.Range("J$r").Value = -0.04
$ws.Range("K$r").Value = ","
$ws.Range("L$r").Value = -0.01
$ws.Range("M$r").Value = ","
$ws.Range("N$r").Value = 0.01
$ws.Range("O$r").Value = ","
$ws.Range("P$r").Value = 0.02
$ws.Range("Q$r").Value = ","
$ws.Range("R$r").Value = 0.01
$ws.Range("S$r").Value = ","
$ws.Range("T$r").Value = 0.03
$ws.Range("U$r").Value = ","
$ws.Range("V$r").Value = 83
$ws.Range("W$r").Value = ", new List<int>() {"
$ws.Range("Y$r").Value = "}, new List<int>(){"
$ws.Range("AA$r").Value = "}"
$ws.Range("AB$r").Value = "),"

# =========================================================================
# 3) Move the active selection to A25 (as recorded in the saved sheet view)
# =========================================================================
$ws.Range("A25").Select()
